$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11").Value = "Giovanni"
$ws.Range("B11").Value = "GDPR"
$ws.Range("C11").Value = 127

$ws.Range("D10").Copy()
$ws.Range("D11").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("D11").Value = Get-Date -Year 2019 -Month 2 -Day 9 -Hour 0 -Minute 0 -Second 0

$ws.Range("B13").Select()
